$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet1 "devTestLogin": the per-row "custom format" flag that every data row
# carried was noise from the original export - clear it off, then drop the
# trailing blank row 11 (dimension shrinks from A1:B11 to A1:B10).
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
for ($i = 2; $i -le 10; $i++) {
    $ws1.Rows.Item($i).ClearFormats()
}
$ws1.Rows.Item(11).Delete()

# ---------------------------------------------------------------------------
# Sheet2 "LoginFunc": nothing to touch directly - once sheet3 becomes the
# active sheet below, Excel naturally clears this sheet's tabSelected flag.
# ---------------------------------------------------------------------------

# ---------------------------------------------------------------------------
# Add the new "userDetails" worksheet as the last tab.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws3.Name = "userDetails"

# Fill the sheet in the same order it was originally authored in, so new
# shared-string entries come out in the same sequence as the saved file
# (header labels, then the two e-mail hyperlinks, then the rest).
$ws3.Cells.Item(1,1).Value = "First  Name "
$ws3.Cells.Item(1,2).Value = "Last Name"
$ws3.Cells.Item(1,3).Value = "Email ID"
$ws3.Cells.Item(1,4).Value = "Phone Number"
$ws3.Cells.Item(1,5).Value = "User Name"
$ws3.Cells.Item(1,6).Value = "Password"
$ws3.Cells.Item(1,7).Value = "Confirm Password"

$ws3.Cells.Item(2,3).Value = "sujith@gmail.com"
$null = $ws3.Hyperlinks.Add($ws3.Range("C2"), "mailto:sujith@gmail.com")
$ws3.Cells.Item(3,3).Value = "sujith@gmail.com"
$null = $ws3.Hyperlinks.Add($ws3.Range("C3"), "mailto:sujith@gmail.com")

$ws3.Cells.Item(1,8).Value = "Login user"
$ws3.Cells.Item(1,9).Value = "Login password"

$ws3.Cells.Item(2,5).Value = "jith"
$ws3.Cells.Item(2,6).Value = "userjith"
$ws3.Cells.Item(2,7).Value = "userjith"
$ws3.Cells.Item(3,5).Value = "jith1"
$ws3.Cells.Item(3,6).Value = "userjith"
$ws3.Cells.Item(3,7).Value = "userjith"

$ws3.Cells.Item(2,1).Value = "Sujith123"
$ws3.Cells.Item(3,1).Value = "Sujith456"
$ws3.Cells.Item(2,2).Value = "css"
$ws3.Cells.Item(3,2).Value = "cs123"

# Numeric phone numbers and the reused "admin"/"useradmin" login columns
$ws3.Cells.Item(2,4).Value = 123456
$ws3.Cells.Item(3,4).Value = 1234566
$ws3.Cells.Item(2,8).Value = "admin"
$ws3.Cells.Item(2,9).Value = "useradmin"
$ws3.Cells.Item(3,8).Value = "admin"
$ws3.Cells.Item(3,9).Value = "useradmin"

# Row 4 - leftover hyperlink-styled placeholder cell, no text
$ws3.Range("C4").Style = $ws3.Range("C2").Style

# Size the columns to fit their content
for ($c = 1; $c -le 9; $c++) {
    $ws3.Columns.Item($c).AutoFit()
}

# Make the new sheet the active tab/selection, matching the saved view state
$ws3.Activate()
$null = $ws3.Range("C6").Select()

Write-Host "edit complete"
